# LOB1252.xlsx content update
#
# Reworks the "Docentes responsaveis" / "Programa resumido" / "Programa" /
# "Avaliacao" block:
#  - the professor name moves up into the "Objetivos" value (row 10)
#  - "Programa resumido" (row 13) gains an explicit A-label and its B/C
#    value becomes "Semestral"
#  - labels in rows 14-21 shift up by one slot
#  - "Programa" (row 15) value becomes the (reused) activation date string
#  - "Avaliação:" (row 17) keeps only its A label; its old B/C value is gone
#  - "Método:" (row 18) value becomes the professor name again
#  - the old Bibliografia text row (old row 22) is removed entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$docente = "5840692 - Diovana Aparecida dos Santos Napoleão"

# --- Row 10 (Objetivos): value replaced with the docente name ---
$ws.Range("B10").Value = $docente
$ws.Range("C10").Value = $docente

# --- Row 13 gains its A label and a brand new value ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14 label becomes "Short syllabus:" (B/C content unchanged) ---
$ws.Range("A14").Value = "Short syllabus:"

# --- Row 15 label becomes "Programa:"; value becomes the reused date text.
#     Copy the whole cell (value+format) from B8/C8, which already hold
#     that exact text as a shared string, so it is not re-parsed as a date
#     and no new number-format style gets minted. ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteAll)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteAll)
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16 label becomes "Syllabus:" (B/C content unchanged) ---
$ws.Range("A16").Value = "Syllabus:"

# --- Row 17 label becomes "Avaliação:"; its old B/C value is removed ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17:C17").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# --- Row 18 label becomes "Método:"; value becomes the docente name again ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = $docente
$ws.Range("C18").Value = $docente
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19 label becomes "Critério:" (B/C content unchanged) ---
$ws.Range("A19").Value = "Critério:"

# --- Row 20 label becomes "Norma de recuperação:" (B/C content unchanged) ---
$ws.Range("A20").Value = "Norma de recuperação:"

# --- Row 21 label becomes "Bibliografia:" (B/C content unchanged) ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120

# --- Old row 22 (old Bibliografia text row) is removed completely ---
$ws.Rows.Item(22).Delete()

Write-Host "done"
